$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New handoff identifiers / timestamps (this commit regenerates the
# handoff package under a new GUID and refreshes the HO xliff dates).
# ---------------------------------------------------------------------
$oldGuid = "cb1d4349-0e19-450f-b29b-f54d6740c1e9"
$newGuid = "3217579c-7ce6-4999-847a-7417e06e9ea5"
$newHash = "3319c4cb329fd23581f4a6b98ca33a48151275a8"

$newFileName   = "$newGuid.md"
$newPathName   = "e2e\$newGuid.md"
$newHoDate     = "2017-02-09 09:34:25"
$newZhName     = "$newGuid.$newHash.zh-cn.xlf"
$newZhHoDate   = "2017-02-09 09:34:04"
$newHbDate     = "0001-01-01 00:00:00"
$newDeName     = "$newGuid.$newHash.de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de0e3e82db7bcddd61cb0e68ce1d15131af63e/e2e/$oldGuid.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newPathName
)

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = $newZhName
$wsZh.Range("H2").Value = $newZhHoDate
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = ""
$wsZh.Range("L2").Value = $newHbDate
$wsZh.Range("Q2").Value = "True"

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de0e3e82db7bcddd61cb0e68ce1d15131af63e/e2e/$oldGuid.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFileName
)

$wsZh.Columns.Item(10).ColumnWidth = 18.6506061553955
$wsZh.Columns.Item(11).ColumnWidth = 21.7054767608643

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = $newDeName
$wsDe.Range("H2").Value = $newHoDate
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = ""
$wsDe.Range("L2").Value = $newHbDate
$wsDe.Range("Q2").Value = "True"

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de0e3e82db7bcddd61cb0e68ce1d15131af63e/e2e/$oldGuid.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFileName
)

$wsDe.Columns.Item(10).ColumnWidth = 18.6506061553955
$wsDe.Columns.Item(11).ColumnWidth = 21.7054767608643
